# "fixed arror in chart" -- the Strassen benchmark numbers in column C of
# Sheet1 were re-measured; update the raw data (which the chart "Strassen's
# Algorithm vs. Transpose Conventional" plots from Sheet1!$C$2:$C$41), then
# leave the selection/active-sheet state the way the author left it when
# they made the fix (on Sheet1, with the corrected column selected).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Corrected "Strassen" timings for Sheet1!C2:C41.
$newC = @(
    1668, 12747, 93643, 93590, 93900, 674354, 669063, 671395, 672281, 667612,
    4748790, 4725772, 4715679, 4799687, 4765063, 4828910, 4788924, 4766968, 4741763, 4731089,
    33662863, 33275755, 33755787, 34315809, 34640412, 35305820, 34554925, 34137404, 34632744, 34763739,
    34131775, 34729742, 34495461, 35604467, 37001781, 35034148, 34622546, 34796305, 34311252, 35468807
)

for ($i = 0; $i -lt $newC.Length; $i++) {
    $ws1.Cells.Item($i + 2, 3).Value = $newC[$i]
}

# Re-point the chart's "Strassen's" series formula at the same range so the
# series picks up the corrected data.
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(2)
$series.Formula = "=SERIES(""Strassen's"",Sheet1!`$A`$2:`$A`$41,Sheet1!`$C`$2:`$C`$41,2)"

# Restore Sheet2's own selection (it's no longer the active tab, but Excel
# still remembers each sheet's last selection independently).
[void]$ws2.Select()
[void]$ws2.Range("H42").Select()

# Finish on Sheet1 with the corrected column selected -- this also makes
# Sheet1 the active tab (activeTab reverts to the default / is omitted).
[void]$ws1.Select()
[void]$ws1.Range("C2:C41").Select()
